# chore: update Sheets via scheduled runner
# Refreshes cached marketboard-derived profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for a handful of leve rows across several job
# sheets. Pure data refresh - no formulas, formatting, or structure changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4805.6816
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4805.6816
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 14417.0448
$ws.Range("N17").Value = -14753.0448

$ws.Range("H31").Value = 2182
$ws.Range("I31").Value = 2182
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6546
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -6316
$ws.Range("N31").ClearContents()

$ws.Range("H64").Value = 4626.5
$ws.Range("I64").Value = 3598.182
$ws.Range("J64").Value = 5654.8184
$ws.Range("K64").Value = 3598.182
$ws.Range("L64").Value = 5654.8184
$ws.Range("M64").Value = -3350.182
$ws.Range("N64").Value = -6150.8184

$ws.Range("H67").Value = 4626.5
$ws.Range("I67").Value = 3598.182
$ws.Range("J67").Value = 5654.8184
$ws.Range("K67").Value = 3598.182
$ws.Range("L67").Value = 5654.8184
$ws.Range("M67").Value = -2740.182
$ws.Range("N67").Value = -7370.8184

$ws.Range("H112").Value = 11364910
$ws.Range("I112").Value = 933.3333
$ws.Range("J112").Value = 12988335
$ws.Range("K112").Value = 2799.9999
$ws.Range("L112").Value = 38965005
$ws.Range("M112").Value = -1691.9999
$ws.Range("N112").Value = -38967221

$ws.Range("H127").Value = 863.5625
$ws.Range("I127").Value = 574.2727
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 1722.8181
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 3237.1819
$ws.Range("N127").Value = -14420

$ws.Range("H129").Value = 1143.1
$ws.Range("I129").Value = 298.5
$ws.Range("J129").Value = 1236.9445
$ws.Range("K129").Value = 895.5
$ws.Range("L129").Value = 3710.8335
$ws.Range("M129").Value = 4104.5
$ws.Range("N129").Value = -13710.8335

$ws.Range("H138").Value = 10484504
$ws.Range("I138").Value = 3969543
$ws.Range("J138").Value = 12823208
$ws.Range("K138").Value = 11908629
$ws.Range("L138").Value = 38469624
$ws.Range("M138").Value = -11903489
$ws.Range("N138").Value = -38479904

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16850.352
$ws.Range("I32").Value = 4207.5625
$ws.Range("J32").Value = 89094.86
$ws.Range("K32").Value = 4207.5625
$ws.Range("L32").Value = 89094.86
$ws.Range("M32").Value = -3920.5625
$ws.Range("N32").Value = -89668.86

$ws.Range("H61").Value = 5721.037
$ws.Range("I61").Value = 5620.56
$ws.Range("J61").Value = 6977
$ws.Range("K61").Value = 5620.56
$ws.Range("L61").Value = 6977
$ws.Range("M61").Value = -5408.56
$ws.Range("N61").Value = -7401

$ws.Range("H133").Value = 60000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 60000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060

$ws.Range("H136").Value = 5721.037
$ws.Range("I136").Value = 5620.56
$ws.Range("J136").Value = 6977
$ws.Range("K136").Value = 16861.68
$ws.Range("L136").Value = 20931
$ws.Range("M136").Value = -14311.68
$ws.Range("N136").Value = -26031

$ws.Range("H139").Value = 47268.125
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 47268.125
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 47268.125
$ws.Range("N139").Value = -57548.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3915.276
$ws.Range("I134").Value = 2579.7778
$ws.Range("J134").Value = 6100.636
$ws.Range("K134").Value = 7739.3334
$ws.Range("L134").Value = 18301.908
$ws.Range("M134").Value = -5204.3334
$ws.Range("N134").Value = -23371.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5366
$ws.Range("I99").Value = 7421.846
$ws.Range("J99").Value = 2936.3635
$ws.Range("K99").Value = 7421.846
$ws.Range("L99").Value = 2936.3635
$ws.Range("M99").Value = -5923.846
$ws.Range("N99").Value = -5932.363499999999

$ws.Range("H126").Value = 5366
$ws.Range("I126").Value = 7421.846
$ws.Range("J126").Value = 2936.3635
$ws.Range("K126").Value = 22265.538
$ws.Range("L126").Value = 8809.0905
$ws.Range("M126").Value = -19795.538
$ws.Range("N126").Value = -13749.0905

$ws.Range("H132").Value = 3725.3572
$ws.Range("I132").Value = 3458.25
$ws.Range("J132").Value = 4081.5
$ws.Range("K132").Value = 10374.75
$ws.Range("L132").Value = 12244.5
$ws.Range("M132").Value = -7844.75
$ws.Range("N132").Value = -17304.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 989.087
$ws.Range("I5").Value = 574.2
$ws.Range("J5").Value = 1767
$ws.Range("K5").Value = 1722.6
$ws.Range("L5").Value = 5301
$ws.Range("M5").Value = -1610.6
$ws.Range("N5").Value = -5525

$ws.Range("H10").Value = 1324.8572
$ws.Range("I10").Value = 89.333336
$ws.Range("J10").Value = 2251.5
$ws.Range("K10").Value = 268.000008
$ws.Range("L10").Value = 6754.5
$ws.Range("M10").Value = -129.000008
$ws.Range("N10").Value = -7032.5

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H98").Value = 512
$ws.Range("I98").Value = 600
$ws.Range("J98").Value = 394.66666
$ws.Range("K98").Value = 1800
$ws.Range("L98").Value = 1183.99998
$ws.Range("M98").Value = -302
$ws.Range("N98").Value = -4179.999980000001

$ws.Range("H109").Value = 3105.4
$ws.Range("I109").Value = 3013.5
$ws.Range("J109").Value = 3166.6667
$ws.Range("K109").Value = 9040.5
$ws.Range("L109").Value = 9500.000100000001
$ws.Range("M109").Value = -8000.5

$ws.Range("H130").Value = 1133.3334
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 1133.3334
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 3400.0002
$ws.Range("N130").Value = -13440.0002
$ws.Range("M130").ClearContents()

$ws.Range("H131").Value = 6174214
$ws.Range("I131").Value = 586.6667
$ws.Range("J131").Value = 6945917.5
$ws.Range("K131").Value = 1760.0001
$ws.Range("L131").Value = 20837752.5
$ws.Range("M131").Value = 3279.9999
$ws.Range("N131").Value = -20847832.5

$ws.Range("H135").Value = 989.087
$ws.Range("I135").Value = 574.2
$ws.Range("J135").Value = 1767
$ws.Range("K135").Value = 5167.8
$ws.Range("L135").Value = 15903
$ws.Range("M135").Value = -2632.8
$ws.Range("N135").Value = -20973

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 28750
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 28750
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 28750
$ws.Range("N4").Value = -28974

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 753.26086
$ws.Range("I22").Value = 548.7857
$ws.Range("J22").Value = 1071.3334
$ws.Range("K22").Value = 548.7857
$ws.Range("L22").Value = 1071.3334
$ws.Range("M22").Value = -253.7857
$ws.Range("N22").Value = -1661.3334

$ws.Range("H27").Value = 753.26086
$ws.Range("I27").Value = 548.7857
$ws.Range("J27").Value = 1071.3334
$ws.Range("K27").Value = 548.7857
$ws.Range("L27").Value = 1071.3334
$ws.Range("M27").Value = -441.7857
$ws.Range("N27").Value = -1285.3334

$ws.Range("H46").Value = 989
$ws.Range("I46").Value = 933.3333
$ws.Range("J46").Value = 1016.8333
$ws.Range("K46").Value = 933.3333
$ws.Range("L46").Value = 1016.8333
$ws.Range("M46").Value = -745.3333
$ws.Range("N46").Value = -1392.8333

$ws.Range("H68").Value = 2260
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2433.3333
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2433.3333
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -3931.3333

$ws.Range("H71").Value = 2260
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2433.3333
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 12166.6665
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -19654.6665

$ws.Range("H132").Value = 5236.2
$ws.Range("I132").Value = 3509.75
$ws.Range("J132").Value = 6829.846
$ws.Range("K132").Value = 10529.25
$ws.Range("L132").Value = 20489.538
$ws.Range("M132").Value = -7999.25
$ws.Range("N132").Value = -25549.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2402.625
$ws.Range("I132").Value = 2534.818
$ws.Range("J132").Value = 1974.9412
$ws.Range("K132").Value = 7604.454000000001
$ws.Range("L132").Value = 5924.8236
$ws.Range("M132").Value = -5074.454000000001
$ws.Range("N132").Value = -10984.8236
